$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column A entirely; remaining columns (old B:F) shift left to A:E,
# keeping their own formatting/values intact.
$ws.Range("A:A").Delete()
